$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the two period/value pairs for GILBERTO ANDRES ZABALETA VILORIA:
# row 16 (period 1705 / 98933) <-> row 18 (period 1707 / 106000)
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 106000
$ws.Range("E18").Value = "1705"
$ws.Range("F18").Value = 98933

# Swap the two periods for CARLOS ARMANDO BUESACO DIAZ GRANADOS:
# row 20 (period 2009) <-> row 21 (period 2010)
$ws.Range("E20").Value = "2010"
$ws.Range("E21").Value = "2009"

$wb.Save()
